$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line "2933 Lamory Pl, Santa Clara CA 95051" into
#    two separate paragraphs: "2933 Lamory Pl" and a new "Santa Clara, CA 95051".
#    Only the first occurrence (the mailing-address block) should change - the
#    "PROPERTY ADDRESS:" table further down keeps the original single-line text,
#    so replace just the first match (wdReplaceOne = 1).
$addrRng = $d.Content
$addrRng.Find.Execute("2933 Lamory Pl, Santa Clara CA 95051", $true, $false, $false, $false, $false,
                       $true, 1, $false, "2933 Lamory Pl", 1)
$endPos = $addrRng.End
$addrRng.InsertParagraphAfter()
$newLineRng = $d.Range($endPos + 1, $endPos + 1)
$newLineRng.InsertAfter("Santa Clara, CA 95051")

# 3. Remove the empty "No Spacing" paragraph that immediately follows
#    "Board of Directors" in the signature block.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -match "Board of Directors") {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
